# "finalized spice scheme, digikey bom, total bom"
#
# Updates the Digi-Key BOM section of the "Total Used Materials" sheet:
#  - header relabel ("Quantity in Package" -> "Quantity Priced",
#    "# to Order" -> "Quantity to Order")
#  - widens the new "Quantity Priced" column (F)
#  - bumps several "Quantity to Order" (H) / package (G) counts
#  - marks several rows as "Already Ordered? = No"
#  - appends a new BOM line (493-14498-ND / 10u electrolytic capacitor)
#  - refreshes the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Quantity / package tweaks across the Digi-Key section ----------------
$ws.Cells.Item(19, 8).Value = 10        # H19  8  -> 10

$ws.Cells.Item(27, 10).Value = "No"     # J27  (new) "No"

$ws.Cells.Item(28, 10).Value = "No"     # J28  (new) "No"

$ws.Cells.Item(29, 8).Value = 20        # H29  2  -> 20
$ws.Cells.Item(29, 10).Value = "No"     # J29  (new) "No"

$ws.Cells.Item(30, 8).Value = 10        # H30  1  -> 10
$ws.Cells.Item(30, 10).Value = "No"     # J30  (new) "No"

$ws.Cells.Item(31, 10).Value = "No"     # J31  (new) "No"

$ws.Cells.Item(32, 8).Value = 50        # H32  (new) 50
$ws.Cells.Item(32, 10).Value = "No"     # J32  (new) "No"

$ws.Cells.Item(33, 8).Value = 10        # H33  1  -> 10
$ws.Cells.Item(33, 10).Value = "No"     # J33  (new) "No"

$ws.Cells.Item(34, 8).Value = 10        # H34  1  -> 10
$ws.Cells.Item(34, 10).Value = "No"     # J34  (new) "No"

$ws.Cells.Item(35, 10).Value = "No"     # J35  (new) "No"

$ws.Cells.Item(36, 10).Value = "No"     # J36  (new) "No"

$ws.Cells.Item(37, 7).Value = 3         # G37  (new) 3
$ws.Cells.Item(37, 8).Value = 2         # H37  1  -> 2
$ws.Cells.Item(37, 10).Value = "No"     # J37  (new) "No"

$ws.Cells.Item(38, 7).Value = 10        # G38  4  -> 10
$ws.Cells.Item(38, 10).Value = "No"     # J38  (new) "No"

# ---- New BOM row: 10u electrolytic capacitor -------------------------------
$ws.Cells.Item(39, 2).Value = "493-14498-ND"
$ws.Cells.Item(39, 3).Value = "PCB"
$ws.Cells.Item(39, 4).Value = "10u electrolytic capacitor"
$ws.Cells.Item(39, 5).Value = 3.98
$ws.Cells.Item(39, 6).Value = 10
$ws.Cells.Item(39, 7).Value = 3
$ws.Cells.Item(39, 8).Value = 2
$ws.Cells.Item(39, 10).Value = "No"

# ---- Header row relabel ---------------------------------------------------
$ws.Cells.Item(1, 6).Value = "Quantity Priced"
$ws.Cells.Item(1, 8).Value = "Quantity to Order"

# ---- Column F ("Quantity Priced") width -----------------------------------
$ws.Columns("F").ColumnWidth = 19

# ---- Refresh active selection ---------------------------------------------
$ws.Range("C33").Select()
